# Bump the "Förändrad" (changed) date in column C by one day (45189 -> 45190,
# i.e. 2023-09-20 -> 2023-09-21) for every data row (rows 2 through 250).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 250
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45190
}
